$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# New matchup rows (Player_1, Points_1, Player_2, Points_2) for week 6 sum input.
$newRows = @(
    @(5, 4, 6, 16),
    @(4, 16, 3, 4),
    @(5, 12, 6, 8),
    @(4, 6, 6, 14),
    @(4, 14, 2, 6),
    @(5, 7, 7, 13),
    @(6, 14, 3, 6),
    @(3, 14, 5, 6),
    @(4, 14, 5, 6),
    @(4, 16, 5, 4),
    @(4, 4, 2, 16),
    @(7, 14, 5, 6),
    @(3, 17, 1, 3),
    @(3, 3, 2, 17),
    @(4, 12, 3, 8),
    @(5, 8, 9, 12)
)

$startRow = 989
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$lastRow = $startRow + $newRows.Count - 1
$selCell = "A" + ($lastRow + 1)
[void]$ws.Range($selCell).Select()

# Scroll the view so the newly-entered rows are visible (mirrors the
# topLeftCell shift from A965 to A986 in the saved workbook view state).
$excel.ActiveWindow.ScrollRow = 986
